$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new quarterly/annual reporting period is inserted as a new first data
# column (before the existing column D). Excel native column-insert semantics
# shift the old D:K data right to E:L and leave a blank column D in place.
$ws.Columns("D").Insert()

# The freshly inserted column D has no formatting yet (Excel gives it the
# left neighbours default style). Copy number-format/style from column E
# (the old column D, now shifted one column right) back onto D so the new
# period lines up visually with the rest of the table (dates, $ amounts, etc).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new reporting periods values in column D.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 35600
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 8800
$ws.Range("D18").Value = 26800
$ws.Range("D20").Value = -17100
$ws.Range("D21").Value = 10700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 9700
$ws.Range("D24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 9200
$ws.Range("D27").Value = 9200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 17100
$ws.Range("D33").Value = 9200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 9200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 9800
$ws.Range("D42").Value = 4700
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 19900
$ws.Range("D49").Value = 19100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1012000
$ws.Range("D57").Value = 800
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 895200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 75800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 116800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 9200
$ws.Range("D83").Value = 1100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 14700
$ws.Range("D91").Value = -400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -29000
$ws.Range("D96").Value = -6200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 16500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 2200
